$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 475

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 63
